# Change the highlight color of the "(1) A (private) thread-safe list of
# "running jobs" ..." bullet from red to cyan (wdTurquoise), including the
# paragraph-mark's own run formatting so the bullet's pilcrow matches the
# rest of the (now cyan-highlighted) list.

$d = $word.ActiveDocument

$find = $d.Content.Find
$find.ClearFormatting()
$find.Text = "A (private) thread-safe list of " + [char]0x201C + "running jobs" + [char]0x201D
$find.Forward = $true
$find.Wrap = 1
$found = $find.Execute()

if ($found) {
    $para = $find.Parent.Paragraphs(1)
    # Setting Font.HighlightColorIndex (rather than Range.HighlightColorIndex)
    # on the full paragraph range also folds the new highlight into the
    # paragraph mark's run properties (<w:pPr><w:rPr>), matching how Word
    # persists highlighting applied to an entire paragraph.
    $para.Range.Font.HighlightColorIndex = 3
    Write-Output "Highlighted 'running jobs' bullet as cyan (wdTurquoise)."
} else {
    Write-Output "WARNING: could not find the 'running jobs' bullet paragraph."
}
